# ==========================================================================
# Edit script reproducing the target diff against before.docx
# ==========================================================================
$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Helper: replace the *entire* paragraph that contains $findText with a
# fresh paragraph built from raw WordprocessingML ($pPrXml / $bodyXml
# are pre-built strings -- NEVER pass two parenthesized expressions in
# a row to a function call in this interpreter, it mis-parses as a
# method call with an empty name). This lets us splice in
# <w:proofErr/> markers (and merge/split runs) exactly like Word's
# proofing pass + editor would.
# --------------------------------------------------------------------
function Replace-ParagraphXml($findText, $pPrXml, $bodyXml) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $rng.Paragraphs(1)
    $prng = $para.Range
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + $pPrXml + $bodyXml + "</w:p>"
    $prng.InsertXML($xml)
}

# --------------------------------------------------------------------
# 1. Title paragraph: insert "q" run + a fresh _GoBack bookmark right
#    before "Yet Another FINAL..."
# --------------------------------------------------------------------
$d.Content.Find.Execute("Yet Another FINAL", $true, $false, $false, $false, $false, $true, 1, $false, "qYet Another FINAL", 2)
$goBackRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --------------------------------------------------------------------
# 2. "Here is data for the ungarbled price (Rs./Kg), garbled price..."
#    -> split so "Rs./" is wrapped by gramStart/gramEnd proofErr marks.
# --------------------------------------------------------------------
$pPr = "<w:pPr><w:pStyle w:val='FirstParagraph'/></w:pPr>"
$body = "<w:r><w:t>Here is data for the ungarbled price (</w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:t>Rs./</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:t>Kg), garbled price (Rs./Kg) and quantity supplied of Cochin peppercorns.</w:t></w:r>"
Replace-ParagraphXml "Here is data for the ungarbled price" $pPr $body

# --------------------------------------------------------------------
# 3. HTMLPreformatted header line: "ungarbled price (Rs./Kg) garbled
#    price (Rs./Kg) quantity supplied" -> same Rs./ proofErr split.
# --------------------------------------------------------------------
$rPr = "<w:rPr><w:rStyle w:val='gd15mcfceub'/><w:rFonts w:ascii='Lucida Console' w:hAnsi='Lucida Console'/><w:u w:val='single'/><w:bdr w:val='none' w:sz='0' w:space='0' w:color='auto' w:frame='1'/></w:rPr>"
$pPr = "<w:pPr><w:pStyle w:val='HTMLPreformatted'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF' w:themeFill='background1'/><w:wordWrap w:val='0'/>" + $rPr + "</w:pPr>"
$body = "<w:r>" + $rPr + "<w:t>ungarbled price (</w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r>" + $rPr + "<w:t>Rs./</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r>" + $rPr + "<w:t>Kg) garbled price (Rs./Kg) quantity supplied</w:t></w:r>"
Replace-ParagraphXml "ungarbled price (Rs./Kg) garbled price" $pPr $body
